# More result table changes.
#
# 1) Switch the active/selected tab back to "Baseline" (it had drifted to
#    "Compression (From LeNet)").
# 2) On "Compression (From LeNet)", move the selection to AD40 and fill in
#    the second ("Relative to base") results block (P29:X41) that had been
#    left blank - mirroring the already-populated first block (C29:K41 /
#    rows 9-18) with seed-sweep utilization numbers, weight counts, and the
#    derived "% of non-zero" / MB figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Compression (From LeNet)")

$utilRow29 = @(0.98949998617172241,0.98949998617172241,0.98949998617172241,0.98949998617172241,0.98949998617172241,0.98949998617172241,0.98949998617172241,0.98949998617172241,0.98949998617172241)
$utilRow30 = @(0.98589998483657837,0.98629999160766602,0.98650002479553223,0.98659998178482056,0.98869997262954712,0.98949998617172241,0.98879998922348022,0.98849999904632568,0.98750001192092896)
$utilRow31 = @(0.99049997329711914,0.99010002613067627,0.98970001935958862,0.98949998617172241,0.98979997634887695,0.98919999599456787,0.98940002918243408,0.98909997940063477,0.98559999465942383)
$utilRow32 = @(0.9878000020980835,0.98750001192092896,0.98729997873306274,0.98739999532699585,0.9869999885559082,0.98619997501373291,0.98559999465942383,0.98559999465942383,0.97890001535415649)
$utilRow33 = @(0.98750001192092896,0.98799997568130493,0.98729997873306274,0.98760002851486206,0.98680001497268677,0.9868999719619751,0.98360002040863037,0.98070001602172852,0.97710001468658447)
$utilRow34 = @(0.98799997568130493,0.9878000020980835,0.98750001192092896,0.98760002851486206,0.98710000514984131,0.98720002174377441,0.98619997501373291,0.98110002279281616,0.97549998760223389)
$utilRow35 = @(0.98530000448226929,0.98500001430511475,0.98540002107620239,0.98519998788833618,0.98489999771118164,0.98400002717971802,0.9836999773979187,0.98379999399185181,0.97189998626708984)
$utilRow36 = @(0.98629999160766602,0.98589998483657837,0.98610001802444458,0.98629999160766602,0.98610001802444458,0.98619997501373291,0.98259997367858887,0.98229998350143433,0.97839999198913574)
$utilRow37 = @(0.98070001602172852,0.98059999942779541,0.98079997301101685,0.98009997606277466,0.98019999265670776,0.97930002212524414,0.9779999852180481,0.97939997911453247,0.97009998559951782)
$utilRow38 = @(0.97469997406005859,0.9747999906539917,0.97460001707077026,0.9747999906539917,0.97390002012252808,0.97369998693466187,0.97219997644424438,0.96649998426437378,0.95770001411437988)

$utilRows = @{
    29 = $utilRow29
    30 = $utilRow30
    31 = $utilRow31
    32 = $utilRow32
    33 = $utilRow33
    34 = $utilRow34
    35 = $utilRow35
    36 = $utilRow36
    37 = $utilRow37
    38 = $utilRow38
}

$cols = @("P","Q","R","S","T","U","V","W","X")

foreach ($r in 29..38) {
    $vals = $utilRows[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $r).Value = $vals[$i]
    }
}

# Row 39: "% of non-zero" sweep values used to build the seed table.
$row39 = @(0.9,0.8,0.7,0.6,0.5,0.4,0.3,0.2,0.1)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "39").Value = $row39[$i]
}

# Row 40: resulting memory size in bytes.
$row40 = @(4534812,4033008,3531204,3029400,2527596,2025792,1523988,1022184,520380)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "40").Value = $row40[$i]
}

# Row 41 already holds "=<cell>40 / 1024^2" formulas; they recalc on their
# own (via Excel's automatic recalculation) once row 40 has values, so they
# are intentionally left untouched here.

# Move this sheet's selection to AD40 (it had been left at K47).
$ws.Range("AD40").Select()

# Re-activate "Baseline" as the selected/visible tab.
$wbBaseline = $wb.Worksheets.Item("Baseline")
$wbBaseline.Activate()
